# Auto-generated edit script applying the Maduin_Profits market-data refresh diff
# across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 9.285714
$ws.Range("I11").Value = 9.285714
$ws.Range("K11").Value = 9.285714
$ws.Range("M11").Value = 130.714286
$ws.Range("H19").Value = 1868.75
$ws.Range("I19").Value = 1712.5
$ws.Range("K19").Value = 1712.5
$ws.Range("M19").Value = -1537.5
$ws.Range("H43").Value = 4157.4
$ws.Range("I43").Value = 4473.25
$ws.Range("J43").Value = 2894
$ws.Range("K43").Value = 4473.25
$ws.Range("L43").Value = 2894
$ws.Range("M43").Value = -4404.25
$ws.Range("N43").Value = -3032
$ws.Range("H53").Value = 142.91667
$ws.Range("I53").Value = 137.4
$ws.Range("K53").Value = 137.4
$ws.Range("M53").Value = 499.6
$ws.Range("H86").Value = 5277.846
$ws.Range("I86").Value = 3099.8
$ws.Range("J86").Value = 6639.125
$ws.Range("K86").Value = 3099.8
$ws.Range("L86").Value = 6639.125
$ws.Range("M86").Value = -1976.8
$ws.Range("N86").Value = -8885.125
$ws.Range("H88").Value = 2976.818
$ws.Range("I88").Value = 2748.75
$ws.Range("J88").Value = 3107.1428
$ws.Range("K88").Value = 2748.75
$ws.Range("L88").Value = 3107.1428
$ws.Range("M88").Value = -2342.75
$ws.Range("N88").Value = -3919.1428
$ws.Range("H89").Value = 5277.846
$ws.Range("I89").Value = 3099.8
$ws.Range("J89").Value = 6639.125
$ws.Range("K89").Value = 15499
$ws.Range("L89").Value = 33195.625
$ws.Range("M89").Value = -9883
$ws.Range("N89").Value = -44427.625
$ws.Range("H91").Value = 2976.818
$ws.Range("I91").Value = 2748.75
$ws.Range("J91").Value = 3107.1428
$ws.Range("K91").Value = 2748.75
$ws.Range("L91").Value = 3107.1428
$ws.Range("M91").Value = -1344.75
$ws.Range("N91").Value = -5915.1428
$ws.Range("H103").Value = 2571.7778
$ws.Range("J103").Value = 4488.1113
$ws.Range("L103").Value = 13464.3339
$ws.Range("N103").Value = -14636.3339

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 9642.857
$ws.Range("I37").Value = 3750
$ws.Range("J37").Value = 17500
$ws.Range("K37").Value = 3750
$ws.Range("L37").Value = 17500
$ws.Range("M37").Value = -3477
$ws.Range("N37").Value = -18046
$ws.Range("H44").Value = 35000
$ws.Range("J44").Value = 35000
$ws.Range("L44").Value = 35000
$ws.Range("N44").Value = -35976
$ws.Range("H55").Value = 25000
$ws.Range("J55").Value = 25000
$ws.Range("L55").Value = 25000
$ws.Range("N55").Value = -25630
$ws.Range("H97").Value = 2382.2727
$ws.Range("I97").Value = 2120.5
$ws.Range("K97").Value = 2120.5
$ws.Range("M97").Value = -1624.5
$ws.Range("H102").Value = 2189.8
$ws.Range("I102").Value = 2189.8
$ws.Range("K102").Value = 2189.8
$ws.Range("M102").Value = -567.8000000000002
$ws.Range("H110").Value = 747.25
$ws.Range("I110").Value = 663
$ws.Range("K110").Value = 663
$ws.Range("M110").Value = 1382
$ws.Range("H139").Value = 95000
$ws.Range("J139").Value = 95000
$ws.Range("L139").Value = 95000
$ws.Range("N139").Value = -105280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3702.6875
$ws.Range("I94").Value = 2906
$ws.Range("K94").Value = 2906
$ws.Range("M94").Value = -2455

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 286.07144
$ws.Range("I22").Value = 333.8889
$ws.Range("J22").Value = 200
$ws.Range("K22").Value = 333.8889
$ws.Range("L22").Value = 200
$ws.Range("M22").Value = 16.11110000000002
$ws.Range("N22").Value = -900
$ws.Range("H99").Value = 5000
$ws.Range("I99").Value = 5000
$ws.Range("K99").Value = 5000
$ws.Range("M99").Value = -3502
$ws.Range("H126").Value = 5000
$ws.Range("I126").Value = 5000
$ws.Range("K126").Value = 15000
$ws.Range("M126").Value = -12530

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 500000260
$ws.Range("I4").Value = 500000260
$ws.Range("K4").Value = 1500000780
$ws.Range("M4").Value = -1500000668
$ws.Range("H23").Value = 189.92308
$ws.Range("I23").Value = 82.5
$ws.Range("J23").Value = 282
$ws.Range("K23").Value = 247.5
$ws.Range("L23").Value = 846
$ws.Range("M23").Value = -12.5
$ws.Range("N23").Value = -1316
$ws.Range("H34").Value = 859.2
$ws.Range("J34").Value = 999.25
$ws.Range("L34").Value = 2997.75
$ws.Range("N34").Value = -3165.75
$ws.Range("H38").Value = 161.55556
$ws.Range("I38").Value = 62.666668
$ws.Range("J38").Value = 359.33334
$ws.Range("K38").Value = 188.000004
$ws.Range("L38").Value = 1078.00002
$ws.Range("M38").Value = 158.999996
$ws.Range("N38").Value = -1772.00002
$ws.Range("H39").Value = 1998.8182
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 1998.8182
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 5996.4546
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -6584.4546
$ws.Range("H55").Value = 2000
$ws.Range("J55").Value = 2000
$ws.Range("L55").Value = 6000
$ws.Range("N55").Value = -6354
$ws.Range("H131").Value = 3633
$ws.Range("J131").Value = 3633
$ws.Range("L131").Value = 10899
$ws.Range("N131").Value = -20979

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 7168.6665
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 7168.6665
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("H97").Value = 4166.6665
$ws.Range("I97").Value = 4000
$ws.Range("K97").Value = 4000
$ws.Range("M97").Value = -3504

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4842.0527
$ws.Range("I46").Value = 3400
$ws.Range("J46").Value = 5357.0713
$ws.Range("K46").Value = 3400
$ws.Range("L46").Value = 5357.0713
$ws.Range("M46").Value = -3212
$ws.Range("N46").Value = -5733.0713
$ws.Range("H61").Value = 900
$ws.Range("I61").Value = 900
$ws.Range("K61").Value = 900
$ws.Range("M61").Value = -698
$ws.Range("H82").Value = 700
$ws.Range("I82").Value = 800
$ws.Range("J82").Value = 650
$ws.Range("K82").Value = 800
$ws.Range("L82").Value = 650
$ws.Range("M82").Value = -439
$ws.Range("N82").Value = -1372
$ws.Range("H85").Value = 700
$ws.Range("I85").Value = 800
$ws.Range("J85").Value = 650
$ws.Range("K85").Value = 800
$ws.Range("L85").Value = 650
$ws.Range("M85").Value = 448
$ws.Range("N85").Value = -3146
$ws.Range("H93").Value = 671.125
$ws.Range("I93").Value = 671.125
$ws.Range("K93").Value = 671.125
$ws.Range("M93").Value = 576.875
$ws.Range("H113").Value = 900
$ws.Range("I113").Value = 900
$ws.Range("K113").Value = 900
$ws.Range("M113").Value = 1270

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1267.4286
$ws.Range("I81").Value = 1267.4286
$ws.Range("K81").Value = 2534.8572
$ws.Range("M81").Value = -1473.8572
$ws.Range("H84").Value = 1267.4286
$ws.Range("I84").Value = 1267.4286
$ws.Range("K84").Value = 12674.286
$ws.Range("M84").Value = -7370.286
$ws.Range("H96").Value = 1906.6428
$ws.Range("I96").Value = 1879.3
$ws.Range("J96").Value = 1975
$ws.Range("K96").Value = 1879.3
$ws.Range("L96").Value = 1975
$ws.Range("M96").Value = -506.3
$ws.Range("N96").Value = -4721
$ws.Range("H100").Value = 6972210
$ws.Range("I100").Value = 11617725
$ws.Range("J100").Value = 3937.5
$ws.Range("K100").Value = 23235450
$ws.Range("L100").Value = 7875
$ws.Range("M100").Value = -23234909
$ws.Range("N100").Value = -8957
$ws.Range("H122").Value = 1199.5
$ws.Range("I122").Value = 1199.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3598.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1148.5
$ws.Range("N122").ClearContents()

